{"js": "const pairs = [\n  [\"26\u00d738=988\", \"29\u00d769=2001\"],\n  [\"97\u00d793=9021\", \"19\u00d723=437\"],\n  [\"66\u00d796=6336\", \"64\u00d749=3136\"],\n  [\"98\u00d747=4606\", \"68\u00d713=884\"],\n  [\"43\u00d726=1118\", \"30\u00d754=1620\"],\n  [\"67\u00d768=4556\", \"93\u00d783=7719\"],\n  [\"73\u00d734=2482\", \"93\u00d783=7719\"],\n  [\"53\u00d740=2120\", \"86\u00d764=5504\"],\n  [\"71\u00d721=1491\", \"49\u00d712=588\"],\n  [\"43\u00d797=4171\", \"83\u00d778=6474\"],\n  [\"90\u00d773=6570\", \"87\u00d738=3306\"],\n  [\"11\u00d734=374\", \"58\u00d743=2494\"],\n  [\"70\u00d782=5740\", \"92\u00d761=5612\"],\n  [\"26\u00d714=364\", \"26\u00d791=2366\"],\n  [\"85\u00d798=8330\", \"65\u00d773=4745\"],\n  [\"71\u00d711=781\", \"59\u00d735=2065\"],\n  [\"42\u00d777=3234\", \"31\u00d728=868\"],\n  [\"56\u00d770=3920\", \"72\u00d756=4032\"],\n  [\"47\u00d743=2021\", \"93\u00d767=6231\"],\n  [\"55\u00d743=2365\", \"42\u00d743=1806\"],\n  [\"94\u00d773=6862\", \"35\u00d716=560\"],\n  [\"42\u00d740=1680\", \"68\u00d784=5712\"],\n  [\"38\u00d794=3572\", \"83\u00d728=2324\"],\n  [\"25\u00d721=525\", \"90\u00d796=8640\"],\n  [\"36\u00d768=2448\", \"68\u00d735=2380\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('26\u00d738=988', '29\u00d769=2001'),\n    @('97\u00d793=9021', '19\u00d723=437'),\n    @('66\u00d796=6336', '64\u00d749=3136'),\n    @('98\u00d747=4606', '68\u00d713=884'),\n    @('43\u00d726=1118', '30\u00d754=1620'),\n    @('67\u00d768=4556', '93\u00d783=7719'),\n    @('73\u00d734=2482', '93\u00d783=7719'),\n    @('53\u00d740=2120', '86\u00d764=5504'),\n    @('71\u00d721=1491', '49\u00d712=588'),\n    @('43\u00d797=4171', '83\u00d778=6474'),\n    @('90\u00d773=6570', '87\u00d738=3306'),\n    @('11\u00d734=374', '58\u00d743=2494'),\n    @('70\u00d782=5740', '92\u00d761=5612'),\n    @('26\u00d714=364', '26\u00d791=2366'),\n    @('85\u00d798=8330', '65\u00d773=4745'),\n    @('71\u00d711=781', '59\u00d735=2065'),\n    @('42\u00d777=3234', '31\u00d728=868'),\n    @('56\u00d770=3920', '72\u00d756=4032'),\n    @('47\u00d743=2021', '93\u00d767=6231'),\n    @('55\u00d743=2365', '42\u00d743=1806'),\n    @('94\u00d773=6862', '35\u00d716=560'),\n    @('42\u00d740=1680', '68\u00d784=5712'),\n    @('38\u00d794=3572', '83\u00d728=2324'),\n    @('25\u00d721=525', '90\u00d796=8640'),\n    @('36\u00d768=2448', '68\u00d735=2380'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $ok = $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 1)\n    if (-not $ok) {\n        throw \"Find/Replace failed for '$oldText' -> '$newText'\"\n    }\n}\n"}
